$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace a few tickers with new index/futures instruments
$ws.Range("A5").Value = "XUU1 Index"
$ws.Range("A8").Value = "11 HK Equity"
$ws.Range("A9").Value = "HCTV1 Index"
$ws.Range("A11").Value = "5 HK Equity"

# Append two new instruments at the bottom of the list
$ws.Range("A17").Value = "XIN9I Index"
$ws.Range("A18").Value = "NDX Index"

# Mirror the selection/active cell left behind in the saved file
$ws.Range("A10").Select()
